$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the _GoBack bookmark that currently sits after the
#    "Supervisor: Phillip Stanley-Marbell" paragraph.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2. Add a new run " 29 May 2019" right after "Date:" (as its own run,
#    not merged into the "Date:" run).
# ---------------------------------------------------------------------
$dateRng = $d.Content
$dateRng.Find.Execute("Date:", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$dateRng.Collapse(0)
$dateRng.InsertAfter(" 29 May 2019")
$newDateRun = $d.Range($dateRng.Start, $dateRng.Start + 12)
# Force the engine to keep this as a distinct run (with its own rPr,
# including the eastAsia font) instead of silently re-merging it with
# the preceding "Date:" run.
$newDateRun.Font.Bold = $true
$newDateRun.Font.Bold = $false

# ---------------------------------------------------------------------
# 3. Split the 29-underscore run that follows "date " into two runs
#    ("_____________" + "________________") and wrap a new _GoBack
#    bookmark around the first of those two runs.
# ---------------------------------------------------------------------
$afterDate = $d.Content
$afterDate.Find.Execute("date ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$afterDate.Collapse(0)

$firstUnderscores = $d.Range($afterDate.Start, $afterDate.Start + 13)
$d.Bookmarks.Add("_GoBack", $firstUnderscores)
